$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 416; this shifts existing rows 416-443 down to 417-444
$ws.Rows.Item(416).Insert()

# Populate the newly inserted row 416 with its data (mirrors the format of surrounding rows)
$ws.Cells.Item(416, 1).Value = 3
$ws.Cells.Item(416, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(416, 3).Value = "Coquimbo"
$ws.Cells.Item(416, 4).Value = 44826
$ws.Cells.Item(416, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(416, 5).Value = 5
$ws.Cells.Item(416, 6).Value = 100112031
$ws.Cells.Item(416, 7).Value = "Poroto verde"
$ws.Cells.Item(416, 8).Value = "Magnum"
$ws.Cells.Item(416, 9).Value = "Primera"
$ws.Cells.Item(416, 10).Value = 68
$ws.Cells.Item(416, 11).Value = 34000
$ws.Cells.Item(416, 12).Value = 35000
$ws.Cells.Item(416, 13).Value = 34559
$ws.Cells.Item(416, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(416, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(416, 16).Value = 1382
$ws.Cells.Item(416, 17).Value = 25
$ws.Cells.Item(416, 18).Value = "Hortaliza"
